# MCH208 collection-level description row.
# The sheet already has a header row (row 1); this adds the single data
# row (row 2) describing the "MCH208" series, matching columns:
#   A=identifier  C=title  D=date_s  E=levelOfDescription
#   F=extentAndMedium  G=notes  (B=alternativeIdentifiers, H=file_path stay blank)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "MCH208"
$ws.Range("C2").Value = "DI SCOTT: GLAUBEN LERNEN IN KRITISCHER ZEIT, SHELL BOYCOTT, DAS LEIDEN BEENDEN, 1987-1989"
$ws.Range("D2").Value = "1987-1989"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: CABINET 1B | GRAP COUNT NUMER: NONE"

# Style the new row (A2 plus C2:H2, skipping the still-blank B2) with the
# same body font family/size used elsewhere in the sheet, in the
# document's dark text theme color.
$newRow = $ws.Range("A2,C2:H2")
foreach ($area in $newRow.Areas) {
  $area.Font.Name = "Calibri"
  $area.Font.Size = 10
  $area.Font.ThemeColor = 1
}

# Re-select the new row and keep the header frozen, as the sheet was left
# after the edit.
$ws.Range("A2:H2").Select()
$excel.ActiveWindow.FreezePanes = $true
